# Update faturamento_diario_lojas.xlsx data:
# Column Y (day 24) values were 0 for all store rows; set them to the
# new reported values, and update column AG (total) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("Y2").Value = 10644.92
$ws.Range("AG2").Value = 254265.14

# Row 3 - Bibi Cell Vieiralves
$ws.Range("Y3").Value = 5142.9
$ws.Range("AG3").Value = 114948.09

# Row 4 - Bibi Cell Manauara
$ws.Range("Y4").Value = 4144
$ws.Range("AG4").Value = 76014.60000000001

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("Y5").Value = 2671
$ws.Range("AG5").Value = 65660.5

# Row 6 - total
$ws.Range("Y6").Value = 22602.82
$ws.Range("AG6").Value = 510888.33
